$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" "329.75"
Set-TextValue "E2" "6.92%"
Set-TextValue "D3" "40.62"
Set-TextValue "E3" "13.34%"
Set-TextValue "D4" "5.948"
Set-TextValue "E4" "16.27%"
Set-TextValue "D5" "0.08142"
Set-TextValue "E5" "5.90%"
Set-TextValue "D6" "4.562"
Set-TextValue "E6" "3.89%"
Set-TextValue "D7" "8.762"
Set-TextValue "E7" "5.59%"
Set-TextValue "D8" "1.954"
Set-TextValue "E8" "5.72%"
Set-TextValue "D9" "2.943"
Set-TextValue "E9" "-0.31%"
Set-TextValue "D10" "0.9435"
Set-TextValue "E10" "2.47%"
Set-TextValue "D11" "0.1313"
Set-TextValue "E11" "17.55%"
Set-TextValue "D12" "0.2009"
Set-TextValue "E12" "8.15%"
Set-TextValue "D13" "0.09252"
Set-TextValue "E13" "5.65%"
Set-TextValue "E14" "2.84%"
Set-TextValue "D15" "0.09631"
Set-TextValue "E15" "1.17%"
Set-TextValue "D16" "0.001321"
Set-TextValue "E16" "-5.29%"
Set-TextValue "D17" "0.006228"
Set-TextValue "E17" "0.04%"
Set-TextValue "D18" "3.374"
Set-TextValue "E18" "0.28%"
Set-TextValue "E19" "1.57%"
Set-TextValue "D20" "7.714"
Set-TextValue "E20" "22.28%"
Set-TextValue "D21" "0.1439"
Set-TextValue "E21" "10.79%"
Set-TextValue "D23" "0.04429"
Set-TextValue "E23" "2.15%"
Set-TextValue "E24" "4.16%"
Set-TextValue "D25" "0.004364"
Set-TextValue "E25" "2.61%"
Set-TextValue "D26" "0.0001189"
Set-TextValue "E26" "-10.75%"
Set-TextValue "D27" "0.0003986"
Set-TextValue "E27" "37.23%"
Set-TextValue "D39" "0.02494"
Set-TextValue "E39" "19.15%"
Set-TextValue "D40" "0.05302"
Set-TextValue "E40" "7.73%"
Set-TextValue "D41" "0.007607"
Set-TextValue "E41" "0.87%"
Set-TextValue "D42" "0.1434"
Set-TextValue "E42" "6.50%"
Set-TextValue "D43" "0.008955"
Set-TextValue "E43" "4.19%"
Set-TextValue "D44" "0.002052"
Set-TextValue "E44" "-1.00%"
Set-TextValue "D45" "0.01057"
Set-TextValue "E45" "25.88%"
Set-TextValue "D46" "0.00006818"
Set-TextValue "E46" "5.43%"
Set-TextValue "E47" "-0.25%"
Set-TextValue "D48" "0.002895"
Set-TextValue "E48" "-12.28%"
Set-TextValue "D49" "0.001798"
Set-TextValue "E49" "24.43%"
Set-TextValue "D50" "0.00002098"
Set-TextValue "E50" "-0.25%"
Set-TextValue "D51" "0.0001998"
Set-TextValue "E51" "-0.25%"
